# Swap the presentation's theme colour scheme from the "Integral" palette
# to the "Office Theme" palette (dk2/lt2/accent1-6/hlink/folHlink) so that
# ppt/theme/theme1.xml ends up carrying the colours that used to live in
# ppt/theme/theme2.xml ("Office Theme"). dk1/lt1 (black/white) and the
# font/format schemes are identical between the two themes, so only the
# colour values need to change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# VBA-style RGB() isn't available in this host, so colours are composed by
# hand as R + G*256 + B*65536 (the standard OLE COLORREF packing that the
# ThemeColorScheme.Item(n).RGB property expects).

# 1 dk1       -> 000000 (unchanged)
$colorScheme.Item(1).RGB = 0x00 + 0x00 * 256 + 0x00 * 65536
# 2 lt1       -> FFFFFF (unchanged)
$colorScheme.Item(2).RGB = 0xFF + 0xFF * 256 + 0xFF * 65536
# 3 dk2       -> 44546A
$colorScheme.Item(3).RGB = 0x44 + 0x54 * 256 + 0x6A * 65536
# 4 lt2       -> E7E6E6
$colorScheme.Item(4).RGB = 0xE7 + 0xE6 * 256 + 0xE6 * 65536
# 5 accent1   -> 5B9BD5
$colorScheme.Item(5).RGB = 0x5B + 0x9B * 256 + 0xD5 * 65536
# 6 accent2   -> ED7D31
$colorScheme.Item(6).RGB = 0xED + 0x7D * 256 + 0x31 * 65536
# 7 accent3   -> A5A5A5
$colorScheme.Item(7).RGB = 0xA5 + 0xA5 * 256 + 0xA5 * 65536
# 8 accent4   -> FFC000
$colorScheme.Item(8).RGB = 0xFF + 0xC0 * 256 + 0x00 * 65536
# 9 accent5   -> 4472C4
$colorScheme.Item(9).RGB = 0x44 + 0x72 * 256 + 0xC4 * 65536
# 10 accent6  -> 70AD47
$colorScheme.Item(10).RGB = 0x70 + 0xAD * 256 + 0x47 * 65536
# 11 hlink    -> 0563C1
$colorScheme.Item(11).RGB = 0x05 + 0x63 * 256 + 0xC1 * 65536
# 12 folHlink -> 954F72
$colorScheme.Item(12).RGB = 0x95 + 0x4F * 256 + 0x72 * 65536
